$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.609405
$ws.Range("H2").Value = 4.828215
$ws.Range("I2").Value = 0.1374279556489526
$ws.Range("J2").Value = 0.1582737093407849
$ws.Range("M2").Value = 0.004501
$ws.Range("N2").Value = 0.009002
$ws.Range("Q2").Value = 0.007243931905
$ws.Range("R2").Value = 0.04346359143
$ws.Range("S2").Value = 0.1374279556489526
$ws.Range("T2").Value = 0.1582737093407849

# Row 3
$ws.Range("I3").Value = 0.4033933174334258
$ws.Range("J3").Value = 0.4645820158786568
$ws.Range("M3").Value = 0.004501
$ws.Range("N3").Value = 0.009002
$ws.Range("Q3").Value = 0.02126316809866666
$ws.Range("R3").Value = 0.127579008592
$ws.Range("S3").Value = 0.4033933174334258
$ws.Range("T3").Value = 0.4645820158786568

# Row 4
$ws.Range("G4").Value = 0.476314
$ws.Range("H4").Value = 1.428942
$ws.Range("I4").Value = 0.04067270778143176
$ws.Range("J4").Value = 0.04684214575631779
$ws.Range("M4").Value = 0.004501
$ws.Range("N4").Value = 0.009002
$ws.Range("Q4").Value = 0.002143889314
$ws.Range("R4").Value = 0.012863335884
$ws.Range("S4").Value = 0.04067270778143176
$ws.Range("T4").Value = 0.04684214575631779

# Row 5
$ws.Range("G5").Value = 4.627222
$ws.Range("H5").Value = 9.254443999999999
$ws.Range("I5").Value = 0.3951209669373822
$ws.Range("J5").Value = 0.3033699161629238
$ws.Range("M5").Value = 0.004501
$ws.Range("N5").Value = 0.009002
$ws.Range("Q5").Value = 0.020827126222
$ws.Range("R5").Value = 0.08330850488799998
$ws.Range("S5").Value = 0.3951209669373822
$ws.Range("T5").Value = 0.3033699161629238

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.27386
$ws.Range("H6").Value = 0.82158
$ws.Range("I6").Value = 0.02338505219880773
$ws.Range("J6").Value = 0.02693221286131667
$ws.Range("M6").Value = 0.004501
$ws.Range("N6").Value = 0.009002
$ws.Range("Q6").Value = 0.00123264386
$ws.Range("R6").Value = 0.007395863159999999
$ws.Range("S6").Value = 0.02338505219880773
$ws.Range("T6").Value = 0.02693221286131667
